# The title paragraph "Inflatie in de supermarkt" (the Title-styled
# paragraph that carries the page-break run) gets its paragraph mark and
# its run switched to the Roboto font (w:rFonts w:ascii="Roboto"
# w:hAnsi="Roboto"), while everything else about the paragraph (style,
# language, text) stays the same.

$d = $word.ActiveDocument

# Locate the exact run of text via Find - this both verifies the target
# text exists and gives us back a Range anchored on that text.
$rng = $d.Content
$found = $rng.Find.Execute("Inflatie in de supermarkt", $true, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find target paragraph text 'Inflatie in de supermarkt'"
}

# Grab the paragraph that owns the found text. Setting Font.Name on the
# paragraph's own Range (which spans the paragraph mark as well as its
# run) updates both the pPr/rPr (paragraph mark formatting) and the
# r/rPr (run formatting) - matching how Word applies a font change when
# the whole paragraph (pilcrow included) is selected.
$para = $rng.Paragraphs(1)
$para.Range.Font.Name = "Roboto"

Write-Output "Applied Roboto font to title paragraph."
